# Generate Report for Handback
# Refresh the timestamps recorded for the "7506985f-1f4c-4eca-ab28-35566073549f"
# file's handoff/handback cycle that is reported on the Overview, zh-cn and
# de-de sheets of the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 7506985f row
$wsOverview.Range("G3").Value = "2016-08-17 08:45:16"

# zh-cn sheet: Correspond Handoff / Handback Datetime for the 7506985f row
$wsZhCn.Range("H3").Value = "2016-08-17 08:44:59"
$wsZhCn.Range("K3").Value = "2016-08-17 08:45:32"

# de-de sheet: Correspond Handoff / Handback Datetime for the 7506985f row
$wsDeDe.Range("H3").Value = "2016-08-17 08:45:16"
$wsDeDe.Range("K3").Value = "2016-08-17 08:45:39"
